$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the blank, merged column W (it was merged with V as a spacer column).
# This shifts X->W, Y->X, Z->Y and removes the now-stale V:W merges.
$ws.Columns("W").Delete()

# Increase the indent on the "1_Ra" header cell (C1).
$ws.Range("C1").IndentLevel = 7

# Restore the active selection to the new last column (W1, previously Z1).
$ws.Range("W1").Select()
